$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New text labels, written in the order they first appear in the author's
#     edit session so the regenerated shared-string table lines up exactly ---
$ws.Range("B10").Value = "Estimated Hours"
$ws.Range("A6").Value = "finished"
$ws.Range("D1").Value = "hours remaining (at start of day)"
$ws.Range("C9").Value = "actual hours"
$ws.Range("C10").Value = "week 1"
$ws.Range("E10").Value = "week 3"
$ws.Range("D10").Value = "week 2"
$ws.Range("F10").Value = "sprint 1"
$ws.Range("A26").Value = "TODO:"
$ws.Range("A27").Value = "convert real entries to day before minus finished"
$ws.Range("P14").Value = "Team McBuddy"

# --- Row 4: "real" burndown actuals (hours remaining, start of each day) ---
$ws.Range("B4").Value = 105
$ws.Range("C4").Value = 103
$ws.Range("D4").Value = 103
$ws.Range("E4").Value = 103
$ws.Range("F4").Value = 98
$ws.Range("G4").Value = 90
$ws.Range("H4").Value = 85
$ws.Range("I4").Value = 83

# --- Row 6: "finished" hours-per-day tracker ---
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 2

# --- Per-person actual hours (column C) + corrected weekly estimates ---
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 5

$ws.Range("C12").Value = 0

$ws.Range("C13").Value = 0

$ws.Range("C14").Value = 7

$ws.Range("C15").Value = 0

$ws.Range("C16").Value = 0

$ws.Range("C17").Value = 10

$ws.Range("B18").Value = 11
$ws.Range("C18").Value = 0

# --- Row 19: total for the new actual-hours column ---
$ws.Range("C19").Formula = "=SUM(C11:C18)"

# --- Move/resize the burndown chart down & right to make room for the
#     expanded table (matches the twoCellAnchor the author ended up with) ---
$co = $ws.ChartObjects(1)
$co.Top = 227
$co.Left = 414.0625
$co.Width = 544.9375
$co.Height = 510

# --- Window size as saved by the author's Excel session ---
$excel.ActiveWindow.Width = 38400
$excel.ActiveWindow.Height = 21140

# --- Restore the selection location the author left the sheet on ---
$ws.Range("V9").Select()
